$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - rename columns to cleaned snake_case names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Fix capitalization in specific municipality/state names (lowercase connector words -> capitalized)
$ws.Range("A8").Value = "Ciudad De México"
$ws.Range("B16").Value = "Taxco De Alarcón"
$ws.Range("B21").Value = "San Miguel El Alto"
$ws.Range("B22").Value = "Tepatitlán De Morelos"
$ws.Range("B40").Value = "Mexquitic De Carmona"
$ws.Range("B46").Value = "Ignacio De La Llave"

# Remove footer rows 54-58 (sample size / source / author / date notes)
$ws.Rows("54:58").Delete()
